$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G (Recorded By) to fit full names
$ws.Columns.Item(7).ColumnWidth = 49.166666666667

# Fill in the "Recorded By" values for column G, rows 2-28 and 30
$ws.Cells.Item(2, 7).Value = 'Dr. Gehan Adel, Dr. Veronia Rafat, Dr. Amira Sobhy, Administrator, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(3, 7).Value = 'Dr. Majorelle Magdy, Dr. Asmaa Reda, Administrator, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Eman Tantawi'
$ws.Cells.Item(4, 7).Value = 'Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(5, 7).Value = 'Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Eman Tantawi'
$ws.Cells.Item(6, 7).Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Mohammad El-Tanany'
$ws.Cells.Item(7, 7).Value = 'Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Lamiaa Ossama'
$ws.Cells.Item(8, 7).Value = 'Dr. Abeer Ragab, Dr. Nada Mohammad'
$ws.Cells.Item(9, 7).Value = 'Dr. Shimaa Ashraf, Dr. Safa Hany'
$ws.Cells.Item(10, 7).Value = 'Dr. Safa Hany'
$ws.Cells.Item(11, 7).Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Cells.Item(12, 7).Value = 'Dr. Madeha Saeed, Dr. Marina Youhanna, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim'
$ws.Cells.Item(13, 7).Value = 'Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa'
$ws.Cells.Item(14, 7).Value = 'Dr. Esraa Samy'
$ws.Cells.Item(15, 7).Value = 'Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef'
$ws.Cells.Item(16, 7).Value = 'Dr. Mohammad Safwat'
$ws.Cells.Item(17, 7).Value = 'Dr. Mohammad Safwat, Dr. Esraa Samy'
$ws.Cells.Item(18, 7).Value = 'Dr. Afnan Fares'
$ws.Cells.Item(19, 7).Value = 'Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Cells.Item(20, 7).Value = 'Dr. Mariam Toma Gerges, Dr. Mohammad Safwat'
$ws.Cells.Item(21, 7).Value = 'Dr. Esraa Samy'
$ws.Cells.Item(22, 7).Value = 'Dr. Alaa Ashraf'
$ws.Cells.Item(23, 7).Value = 'Menna tuâ€™Allah Gamil'
$ws.Cells.Item(24, 7).Value = 'Dr. Youstina Gamil, Dr. Sarah Mahdy'
$ws.Cells.Item(25, 7).Value = 'Menna tuâ€™Allah Gamil, Dr. Nouran Mahmoud'
$ws.Cells.Item(26, 7).Value = 'Dr. Nancy Abd Al-Shafy'
$ws.Cells.Item(27, 7).Value = 'Dr. Hana Amr, Dr. Nourham Mostafa'
$ws.Cells.Item(28, 7).Value = 'Dr. Aya Emad, Dr. Maryam Ashraf'
$ws.Cells.Item(30, 7).Value = 'Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Shorok Mohammad'
